$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 164 - 四方坪站充电量(kw) - date 2025-11-21 (serial 45982)
$row164 = @{
    A = 45982
    B = "四方坪站充电量(kw)"
    C = 536.3889999999999
    D = 900.40700000000004
    E = 610.46399999999994
    F = 322.66899999999998
    G = 219.5
    H = 786.43300000000011
    I = 511.53200000000004
    J = 207.00299999999999
    K = 144.37299999999999
    L = 162.679
    M = 97.82
    N = 172.48
    O = 688.69999999999993
    P = 1665.2529999999992
    Q = 639.9409999999998
    R = 391.43099999999998
    S = 306.68900000000002
    T = 370.34899999999999
    U = 112.849
    V = 76.63
    W = 128.53
    X = 223.54099999999997
    Y = 94.99
    Z = 30.22
}

# Row 165 - 高岭站充电量(kw) - date 2025-11-21 (serial 45982)
$row165 = @{
    A = 45982
    B = "高岭站充电量(kw)"
    C = 702.70600000000013
    D = 251.00300000000001
    E = 170.58700000000002
    F = 37.859000000000002
    G = 94.738
    H = 348.56100000000004
    I = 175.81700000000001
    J = 119.809
    K = 205.34500000000003
    L = 222.06299999999999
    M = 175.80500000000001
    N = 298.74899999999997
    O = 337.26600000000002
    P = 436.37900000000008
    Q = 389.06299999999999
    R = 435.89299999999992
    S = 115.999
    T = 126.24599999999998
    U = 36.57
    V = 77.968999999999994
    W = 83.352000000000004
    X = 72.13300000000001
    Y = 34.823
    Z = 29.445
}

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

foreach ($col in $columns) {
    $addr164 = "{0}164" -f $col
    $ws.Range($addr164).Value = $row164[$col]
}
foreach ($col in $columns) {
    $addr165 = "{0}165" -f $col
    $ws.Range($addr165).Value = $row165[$col]
}

$ws.Range("C168").Select()
